$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties columns, matching the header
# formatting (bold, centered, bordered) already used by the other
# header cells in row 1.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-57: season record values (Wins=79, Losses=83, Ties=0)
# for every player row in the sheet.
for ($row = 2; $row -le 57; $row++) {
    $ws.Cells.Item($row, 30).Value = 79
    $ws.Cells.Item($row, 31).Value = 83
    $ws.Cells.Item($row, 32).Value = 0
}
